$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Title / byline text, appears twice (Heading1 and bold byline near the end)
Replace-Text "Play Irish Magic for Free & Win Big - Review" "Play Irish Magic for Free - Review & Gameplay Features"

# "What we like" bullet list
Replace-Text "Leprechaun Wild and Free Spins mode" "Wide range of bet limits"
Replace-Text "High payout Progressive Jackpot" "Progressive Jackpot feature"
Replace-Text "Generous prizes in base game" "Leprechaun Wild and Free Spins features"
Replace-Text "Dynamic soundtrack and design" "Visually and audibly appealing design"

# "What we don't like" bullet list
Replace-Text "High volatility" "High volatility may result in consecutive failed spins"
Replace-Text "Outdated graphics" "Graphics have a 90s-style look"

# Meta description (italic)
Replace-Text "Play Irish Magic for free and enjoy Leprechaun Wild symbol, progressive jackpot, dynamic soundtrack, and generous prizes in base game. Our review explains." "Play Irish Magic, a visually appealing slot game with a progressive jackpot and exciting features."
